# Fill in the newly-added skills (IDs 8-13) in the "Skills" sheet's second
# table (Mage skills), and refresh the Name column's auto-fit width plus the
# current selection to reflect where editing left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skills")

$ws.Range("B14").Value = "Firewall"
$ws.Range("B15").Value = "Hailstone"
$ws.Range("B16").Value = "Power-sphere"
$ws.Range("B17").Value = "Lightning"
$ws.Range("B18").Value = "Levitation"
$ws.Range("B19").Value = "Armageddon"

# Column B ("Name") is best-fit; recompute its width now that longer names
# (e.g. "Power-sphere") have been entered.
$ws.Columns.Item(2).AutoFit()

# Leave the selection where work stopped, on the next empty row's Name cell.
[void]$ws.Range("B20").Select()
